# Updates cryptos list values (price/volume columns D & E) to match
# the latest scrape, and fixes the EnergySwap/PaxDollar row swap
# (rows 47-48 had their data in the wrong order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds numeric-looking strings (e.g. "306.76", "1.001")
# that must stay as plain text (as in the source data) rather than being
# auto-coerced into floating point numbers by Excel's input parser, which
# would introduce binary-rounding artifacts (e.g. 306.76 -> 306.75999...).
# Forcing a text number format on the column before writing keeps every
# assignment below a literal text value.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.190.51'
$ws.Range("D3").Value = '1.902.62'
$ws.Range("E3").Value = '  +1.46%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = '306.76'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").Value = '0.5235'
$ws.Range("E7").Value = '  +1.52%  '
$ws.Range("D8").Value = '0.3769'
$ws.Range("E8").Value = '  +1.05%  '
$ws.Range("D9").Value = '0.07245'
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("D10").Value = '21.17'
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("D11").Value = '0.8989'
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").Value = '0.08412'
$ws.Range("E12").Value = '  +11.30%  '
$ws.Range("D13").Value = '1.908.12'
$ws.Range("E13").Value = '  +1.85%  '
$ws.Range("D14").Value = '94.57'
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("D15").Value = '5.268'
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("E17").Value = '  +0.85%  '
$ws.Range("E18").Value = '  +1.57%  '
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").Value = '27.228.47'
$ws.Range("E20").Value = '  +0.96%  '
$ws.Range("D21").Value = '5.054'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").Value = '2.142.08'
$ws.Range("E22").Value = '  +3.01%  '
$ws.Range("E23").Value = '  +1.70%  '
$ws.Range("D24").Value = '6.417'
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").Value = '2.286'
$ws.Range("E25").Value = '  +8.08%  '
$ws.Range("D26").Value = '146.66'
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("E28").Value = '  +0.67%  '
$ws.Range("D29").Value = '114.84'
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").Value = '4.927'
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("D31").Value = '4.787'
$ws.Range("E31").Value = '  +0.46%  '
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").Value = '0.8147'
$ws.Range("E33").Value = '  +8.21%  '
$ws.Range("D34").Value = '0.05056'
$ws.Range("E34").Value = '  +0.49%  '
$ws.Range("D35").Value = '1.239'
$ws.Range("E35").Value = '  +5.67%  '
$ws.Range("D36").Value = '2.949'
$ws.Range("E36").Value = '  -1.40%  '
$ws.Range("D37").Value = '3.381'
$ws.Range("E37").Value = '  +3.65%  '
$ws.Range("D38").Value = '2.563'
$ws.Range("E38").Value = '  +2.88%  '
$ws.Range("D39").Value = '0.5698'
$ws.Range("E39").Value = '  +2.12%  '
$ws.Range("D40").Value = '0.01976'
$ws.Range("E40").Value = '  -0.88%  '
$ws.Range("D41").Value = '1.072'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = '6.644'
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("D43").Value = '8.958'
$ws.Range("E43").Value = '  +2.76%  '
$ws.Range("D44").Value = '118.25'
$ws.Range("E44").Value = '  +2.18%  '
$ws.Range("D45").Value = '0.1511'
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("D46").Value = '0.4825'
$ws.Range("E46").Value = '  +1.07%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '10.13'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").Value = '1.613'
$ws.Range("E49").Value = '  +3.08%  '
$ws.Range("D50").Value = '37.41'
$ws.Range("E50").Value = '  +0.86%  '
$ws.Range("D51").Value = '63.64'
$ws.Range("E51").Value = '  +0.27%  '
